# Update view counts (column F) across the sheets to reflect the latest
# generated output (commit: "Update gh-pages to output generated at 456a3b4").
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local life)
# Sheet 4 = 全部类型 (All types - combined view of sheets 1-3)

$wb = $excel.ActiveWorkbook

$wsExhibition  = $wb.Worksheets.Item(1)
$wsPerformance = $wb.Worksheets.Item(2)
$wsLocalLife   = $wb.Worksheets.Item(3)
$wsAll         = $wb.Worksheets.Item(4)

# 展览 (Exhibitions)
$wsExhibition.Range("F4").Value  = 20621
$wsExhibition.Range("F6").Value  = 2637
$wsExhibition.Range("F7").Value  = 785
$wsExhibition.Range("F18").Value = 4
$wsExhibition.Range("F20").Value = 16

# 演出 (Performances)
$wsPerformance.Range("F22").Value = 38

# 本地生活 (Local life)
$wsLocalLife.Range("F4").Value = 659
$wsLocalLife.Range("F5").Value = 1463

# 全部类型 (All types)
$wsAll.Range("F4").Value  = 659
$wsAll.Range("F5").Value  = 1463
$wsAll.Range("F8").Value  = 20621
$wsAll.Range("F14").Value = 2637
$wsAll.Range("F15").Value = 785
$wsAll.Range("F37").Value = 4
$wsAll.Range("F40").Value = 16
$wsAll.Range("F49").Value = 38
